$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row 168 (data row, run_id=167)
$ws.Range("A168").Value = 167
$ws.Range("B168").Value = 1
$ws.Range("C168").Value = "2024-06-18 12:23:55"
$ws.Range("D168").Value = 200
$ws.Range("E168").Value = 21

# New row 169 (data row, run_id=168)
$ws.Range("A169").Value = 168
$ws.Range("B169").Value = 2
$ws.Range("C169").Value = "2024-06-18 12:23:55"
$ws.Range("D169").Value = 200
$ws.Range("E169").Value = 0
